$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> [D value, E value] (column D = 4, column E = 5).
# A leading apostrophe marks values that look like plain numbers (e.g.
# "586.65") so Excel keeps them as literal text instead of silently
# re-interpreting them as numeric cells - exactly what typing '586.65
# into a cell does.
$updates = @{
    2  = @("64.289.02", "  +0.76%  ")
    3  = @("3.485.09", "  +0.50%  ")
    4  = @($null, "  -0.22%  ")
    5  = @("'586.65", "  +0.88%  ")
    6  = @("'133.88", "  +2.58%  ")
    7  = @("3.485.69", "  +0.54%  ")
    8  = @($null, "  -0.05%  ")
    9  = @($null, "  -0.62%  ")
    10 = @($null, "  +0.36%  ")
    11 = @("'7.19", "  +0.82%  ")
    12 = @("'0.375", "  -2.03%  ")
    13 = @("4.081.21", "  +0.38%  ")
    14 = @($null, "  +2.02%  ")
    15 = @($null, "  +1.69%  ")
    16 = @("3.484.93", "  +0.15%  ")
    17 = @("64.345.83", "  +1.42%  ")
    18 = @("'25.21", "  -8.64%  ")
    19 = @("'9.98", "  +1.37%  ")
    20 = @("'5.65", "  +0.69%  ")
    21 = @("'13.71", "  -3.23%  ")
    22 = @("'384.79", "  -1.25%  ")
    23 = @("'0.566", "  -1.04%  ")
    24 = @("3.626.79", "  +0.21%  ")
    25 = @("'74.08", "  +1.45%  ")
    26 = @($null, "  +0.01%  ")
    27 = @("'5.70", "  -0.44%  ")
    28 = @($null, "  +5.24%  ")
    29 = @("'1.55", "  +0.10%  ")
    30 = @("'0.999", "  +0.71%  ")
    31 = @("'7.44", "  +0.14%  ")
    32 = @($null, "  -0.27%  ")
    33 = @("'8.20", "  +1.08%  ")
    34 = @("3.509.22", "  +0.96%  ")
    35 = @($null, "  -0.01%  ")
    36 = @($null, "  +2.07%  ")
    37 = @("'23.35", "  -1.39%  ")
    38 = @($null, "  +1.88%  ")
    39 = @($null, "  -1.28%  ")
    40 = @($null, "  -1.69%  ")
    41 = @("'162.28", "  -4.22%  ")
    42 = @("'0.0777", "  -2.70%  ")
    43 = @("'0.801", "  -0.60%  ")
    44 = @("'25.73", "  +1.63%  ")
    45 = @($null, "  -0.15%  ")
    46 = @($null, "  +0.62%  ")
    47 = @("'4.38", "  +1.34%  ")
    48 = @($null, "  +0.61%  ")
    49 = @("'1.65", "  +2.20%  ")
    50 = @("2.464.32", "  +2.07%  ")
    51 = @("'6.73", "  -1.38%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($dVal -ne $null) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    if ($eVal -ne $null) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}
